$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.156.08"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").Value = "2.505.20"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.60"
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.67"
$ws.Range("E6").Value = "  -1.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.525"
$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  -2.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.59"
$ws.Range("E10").Value = "  -2.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.16"
$ws.Range("E11").Value = "  +7.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0811"
$ws.Range("E12").Value = "  -0.26%  "

$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.10"
$ws.Range("E14").Value = "  -1.64%  "

$ws.Range("D15").Value = "2.895.00"
$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("D16").Value = "2.499.00"
$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.836"
$ws.Range("E17").Value = "  -1.90%  "

$ws.Range("D18").Value = "48.007.99"
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("E19").Value = "  -1.97%  "

$ws.Range("E20").Value = "  +1.24%  "

$ws.Range("D21").Value = "0.0₃0939"
$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.77"
$ws.Range("E22").Value = "  -0.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "277.46"
$ws.Range("E23").Value = "  +12.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.44"
$ws.Range("E24").Value = "  +0.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  -0.61%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.87"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("E28").Value = "  -2.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.140"
$ws.Range("E29").Value = "  +0.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.32"
$ws.Range("E30").Value = "  +0.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.10"
$ws.Range("E31").Value = "  +0.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.53"
$ws.Range("E32").Value = "  -0.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.53"
$ws.Range("E33").Value = "  -2.44%  "

$ws.Range("E34").Value = "  -0.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"
$ws.Range("E35").Value = "  -0.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0783"
$ws.Range("E36").Value = "  -1.07%  "

$ws.Range("E37").Value = "  -0.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.63"
$ws.Range("E38").Value = "  -0.90%  "

$ws.Range("E39").Value = "  -2.73%  "

$ws.Range("E40").Value = "  -0.61%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.21"
$ws.Range("E41").Value = "  +1.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.43"
$ws.Range("E43").Value = "  -4.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0301"
$ws.Range("E44").Value = "  +0.48%  "

$ws.Range("D45").Value = "2.019.57"
$ws.Range("E45").Value = "  +0.87%  "

$ws.Range("E46").Value = "  +2.82%  "

$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("E48").Value = "  +1.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.00"
$ws.Range("E49").Value = "  -0.43%  "

$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.13"
$ws.Range("E51").Value = "  +3.14%  "
